$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 45740.01041666666, 1751)
    ,@(3, 45740.02083333334, 1749)
    ,@(4, 45740.03125, 1744)
    ,@(5, 45740.04166666666, 1730)
    ,@(6, 45740.05208333334, 1536)
    ,@(7, 45740.0625, 1531)
    ,@(8, 45740.07291666666, 1527)
    ,@(9, 45740.08333333334, 1519)
    ,@(10, 45740.09375, 1305)
    ,@(11, 45740.10416666666, 1301)
    ,@(12, 45740.11458333334, 1295)
    ,@(13, 45740.125, 1285)
    ,@(14, 45740.13541666666, 1100)
    ,@(15, 45740.14583333334, 1096)
    ,@(16, 45740.15625, 1079)
    ,@(17, 45740.16666666666, 1064)
    ,@(18, 45740.17708333334, 896)
    ,@(19, 45740.1875, 890)
    ,@(20, 45740.19791666666, 888)
    ,@(21, 45740.20833333334, 874)
    ,@(22, 45740.21875, 814)
    ,@(23, 45740.22916666666, 811)
    ,@(24, 45740.23958333334, 816)
    ,@(25, 45740.25, 811)
    ,@(26, 45740.26041666666, 667)
    ,@(27, 45740.27083333334, 665)
    ,@(28, 45740.28125, 667)
    ,@(29, 45740.29166666666, 662)
    ,@(30, 45740.30208333334, 531)
    ,@(31, 45740.3125, 527)
    ,@(32, 45740.32291666666, 525)
    ,@(33, 45740.33333333334, 519)
    ,@(34, 45740.34375, 408)
    ,@(35, 45740.35416666666, 404)
    ,@(36, 45740.36458333334, 402)
    ,@(37, 45740.375, 401)
    ,@(38, 45740.38541666666, 296)
    ,@(39, 45740.39583333334, 296)
    ,@(40, 45740.40625, 297)
    ,@(41, 45740.41666666666, 297)
    ,@(42, 45740.42708333334, 283)
    ,@(43, 45740.4375, 284)
    ,@(44, 45740.44791666666, 285)
    ,@(45, 45740.45833333334, 286)
    ,@(46, 45740.46875, 292)
    ,@(47, 45740.47916666666, 293)
    ,@(48, 45740.48958333334, 295)
    ,@(49, 45740.5, 296)
    ,@(50, 45740.51041666666, 314)
    ,@(51, 45740.52083333334, 316)
    ,@(52, 45740.53125, 318)
    ,@(53, 45740.54166666666, 320)
    ,@(54, 45740.55208333334, 355)
    ,@(55, 45740.5625, 357)
    ,@(56, 45740.57291666666, 360)
    ,@(57, 45740.58333333334, 363)
    ,@(58, 45740.59375, 409)
    ,@(59, 45740.60416666666, 411)
    ,@(60, 45740.61458333334, 414)
    ,@(61, 45740.625, 417)
    ,@(62, 45740.63541666666, 477)
    ,@(63, 45740.64583333334, 480)
    ,@(64, 45740.65625, 483)
    ,@(65, 45740.66666666666, 486)
    ,@(66, 45740.67708333334, 615)
    ,@(67, 45740.6875, 619)
    ,@(68, 45740.69791666666, 624)
    ,@(69, 45740.70833333334, 629)
    ,@(70, 45740.71875, 808)
    ,@(71, 45740.72916666666, 814)
    ,@(72, 45740.73958333334, 821)
    ,@(73, 45740.75, 827)
    ,@(74, 45740.76041666666, 957)
    ,@(75, 45740.77083333334, 960)
    ,@(76, 45740.78125, 963)
    ,@(77, 45740.79166666666, 965)
    ,@(78, 45740.80208333334, 980)
    ,@(79, 45740.8125, 981)
    ,@(80, 45740.82291666666, 981)
    ,@(81, 45740.83333333334, 982)
    ,@(82, 45740.84375, 911)
    ,@(83, 45740.85416666666, 911)
    ,@(84, 45740.86458333334, 910)
    ,@(85, 45740.875, 910)
    ,@(86, 45740.88541666666, 843)
    ,@(87, 45740.89583333334, 842)
    ,@(88, 45740.90625, 841)
    ,@(89, 45740.91666666666, 840)
    ,@(90, 45740.92708333334, 788)
    ,@(91, 45740.9375, 785)
    ,@(92, 45740.94791666666, 783)
    ,@(93, 45740.95833333334, 780)
    ,@(94, 45740.96875, 0)
    ,@(95, 45740.97916666666, 0)
    ,@(96, 45740.98958333334, 0)
    ,@(97, 45741, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $a = $row[1]
    $b = $row[2]
    $ws.Cells.Item($r, 1).Value2 = $a
    $ws.Cells.Item($r, 2).Value2 = $b
}
